{"js": "// Flyer template text update (v1.20):\n// \"\u0420\u0435\u043c\u043e\u043d\u0442 \u0434\u043e\u043c\u0430\" -> \"\u0422\u0435\u043a. \u0440\u0435\u043c\u043e\u043d\u0442\" in the services table.\nconst searchResults = context.document.body.search(\"\u0420\u0435\u043c\u043e\u043d\u0442 \u0434\u043e\u043c\u0430\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"\u0422\u0435\u043a. \u0440\u0435\u043c\u043e\u043d\u0442\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Flyer template text update (v1.20):\n# \"\u0420\u0435\u043c\u043e\u043d\u0442 \u0434\u043e\u043c\u0430\" -> \"\u0422\u0435\u043a. \u0440\u0435\u043c\u043e\u043d\u0442\" in the services table.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.Execute(\"\u0420\u0435\u043c\u043e\u043d\u0442 \u0434\u043e\u043c\u0430\", $false, $false, $false, $false, $false, $true, 1, $false, \"\u0422\u0435\u043a. \u0440\u0435\u043c\u043e\u043d\u0442\", 2) | Out-Null\n"}
